$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 16 (new): "Excel动手实验室 - 数据排序和筛选"  (Title and Content layout)
# ---------------------------------------------------------------------------
$s16 = $p.Slides.Add(16, 2)

$s16.Shapes.Item(1).TextFrame.TextRange.Text = "Excel动手实验室 - 数据排序和筛选"

$body16 = $s16.Shapes.Item(2).TextFrame
$body16.TextRange.Text = "排序 & 筛选简介`r排序`r单列排序`r多列排序`r自定义排序`r按照颜色排序`r筛选`r自动筛选`r自定义筛选"
$body16.TextRange.Paragraphs(3).IndentLevel = 2
$body16.TextRange.Paragraphs(4).IndentLevel = 2
$body16.TextRange.Paragraphs(5).IndentLevel = 2
$body16.TextRange.Paragraphs(6).IndentLevel = 2
$body16.TextRange.Paragraphs(8).IndentLevel = 2
$body16.TextRange.Paragraphs(9).IndentLevel = 2
$body16.AutoSize = 2

# ---------------------------------------------------------------------------
# Slide 17 (new): "Excel动手实验室 – 高级筛选"  (Title and Content layout)
# ---------------------------------------------------------------------------
$s17 = $p.Slides.Add(17, 2)

$s17.Shapes.Item(1).TextFrame.TextRange.Text = "Excel动手实验室 – 高级筛选"

$body17 = $s17.Shapes.Item(2).TextFrame
$body17.TextRange.Text = "高级筛选`r重复项`r条件区域`r条件区域自定义条件`r复制筛选结果"
$body17.TextRange.Paragraphs(2).IndentLevel = 2
$body17.TextRange.Paragraphs(3).IndentLevel = 2
$body17.TextRange.Paragraphs(4).IndentLevel = 2
$body17.TextRange.Paragraphs(5).IndentLevel = 2

# ---------------------------------------------------------------------------
# Slide 2: mark the existing "兼容性" paragraph run clean (dirty="0")
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$run = $slide2.Shapes.Item(2).TextFrame.TextRange.Paragraphs(11).Runs(1)
$run.Text = $run.Text
